# Updates the Overview sheet's quarterly dollar-cumulative income-statement
# table: drop the oldest ("6 mahe montahi be 1399/06") quarter column,
# shift every remaining quarter one column to the left (D<-E, E<-F, ... L<-M),
# and append the newly published quarter ("12 mahe montahi be 1401/12",
# published 1402-01-29) into column M - for the period headers (row 8),
# publish dates (row 9) and every financial-statement data row (11-26).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 11-26 quarterly data: shift one quarter left, append newest quarter ---
# Row 11
$ws.Range("D11").Value = 711004
$ws.Range("E11").Value = 1220478
$ws.Range("F11").Value = 418187
$ws.Range("G11").Value = 867404
$ws.Range("H11").Value = 1290101
$ws.Range("I11").Value = 1728181
$ws.Range("J11").Value = 422195
$ws.Range("K11").Value = 684575
$ws.Range("L11").Value = 972894
$ws.Range("M11").Value = 1163255

# Row 12
$ws.Range("D12").Value = -245509
$ws.Range("E12").Value = -410832
$ws.Range("F12").Value = -137174
$ws.Range("G12").Value = -288772
$ws.Range("H12").Value = -514480
$ws.Range("I12").Value = -630616
$ws.Range("J12").Value = -176732
$ws.Range("K12").Value = -351034
$ws.Range("L12").Value = -539410
$ws.Range("M12").Value = -620833

# Row 13
$ws.Range("D13").Value = 465496
$ws.Range("E13").Value = 743965
$ws.Range("F13").Value = 281013
$ws.Range("G13").Value = 578633
$ws.Range("H13").Value = 775621
$ws.Range("I13").Value = 1097565
$ws.Range("J13").Value = 245463
$ws.Range("K13").Value = 333541
$ws.Range("L13").Value = 433485
$ws.Range("M13").Value = 542422

# Row 14
$ws.Range("D14").Value = -2355
$ws.Range("E14").Value = -14033
$ws.Range("F14").Value = -1281
$ws.Range("G14").Value = -13186
$ws.Range("H14").Value = -18733
$ws.Range("I14").Value = -28704
$ws.Range("J14").Value = -6613
$ws.Range("K14").Value = -13502
$ws.Range("L14").Value = -19980
$ws.Range("M14").Value = -28196

# Row 17
$ws.Range("D17").Value = 463140
$ws.Range("E17").Value = 729933
$ws.Range("F17").Value = 279732
$ws.Range("G17").Value = 565447
$ws.Range("H17").Value = 756888
$ws.Range("I17").Value = 1068862
$ws.Range("J17").Value = 238850
$ws.Range("K17").Value = 320039
$ws.Range("L17").Value = 413504
$ws.Range("M17").Value = 514225

# Row 18
$ws.Range("D18").Value = -272
$ws.Range("E18").Value = -265
$ws.Range("F18").Value = "-"
$ws.Range("G18").Value = "-"
$ws.Range("H18").Value = "-"
$ws.Range("I18").Value = "-"
$ws.Range("J18").Value = "-"
$ws.Range("K18").Value = "-"
$ws.Range("L18").Value = "-"
$ws.Range("M18").Value = "-"

# Row 19
$ws.Range("D19").Value = 7971
$ws.Range("E19").Value = 24595
$ws.Range("F19").Value = 11058
$ws.Range("G19").Value = 20079
$ws.Range("H19").Value = 31773
$ws.Range("I19").Value = 42376
$ws.Range("J19").Value = 12272
$ws.Range("K19").Value = 26480
$ws.Range("L19").Value = 55789
$ws.Range("M19").Value = 65182

# Row 20
$ws.Range("D20").Value = 470839
$ws.Range("E20").Value = 819944
$ws.Range("F20").Value = 290790
$ws.Range("G20").Value = 585526
$ws.Range("H20").Value = 788661
$ws.Range("I20").Value = 1111238
$ws.Range("J20").Value = 251123
$ws.Range("K20").Value = 346519
$ws.Range("L20").Value = 469293
$ws.Range("M20").Value = 579407

# Row 21
$ws.Range("D21").Value = -41087
$ws.Range("E21").Value = -63947
$ws.Range("F21").Value = -33055
$ws.Range("G21").Value = -65166
$ws.Range("H21").Value = -87671
$ws.Range("I21").Value = -89924
$ws.Range("J21").Value = -16515
$ws.Range("K21").Value = -28427
$ws.Range("L21").Value = -28565
$ws.Range("M21").Value = "-"

# Row 22
$ws.Range("D22").Value = 429752
$ws.Range("E22").Value = 755997
$ws.Range("F22").Value = 257735
$ws.Range("G22").Value = 520360
$ws.Range("H22").Value = 700990
$ws.Range("I22").Value = 1021313
$ws.Range("J22").Value = 234608
$ws.Range("K22").Value = 318092
$ws.Range("L22").Value = 440728
$ws.Range("M22").Value = 579407

# Row 24
$ws.Range("D24").Value = 429752
$ws.Range("E24").Value = 755997
$ws.Range("F24").Value = 257735
$ws.Range("G24").Value = 520360
$ws.Range("H24").Value = 700990
$ws.Range("I24").Value = 1021313
$ws.Range("J24").Value = 234608
$ws.Range("K24").Value = 318092
$ws.Range("L24").Value = 440728
$ws.Range("M24").Value = 579407

# Row 26
$ws.Range("D26").Value = 252317
$ws.Range("E26").Value = 245464
$ws.Range("F26").Value = 237892
$ws.Range("G26").Value = 223965
$ws.Range("H26").Value = 457438
$ws.Range("I26").Value = 451008
$ws.Range("J26").Value = 404972
$ws.Range("K26").Value = 394109
$ws.Range("L26").Value = 374220
$ws.Range("M26").Value = 545495

# --- Row 8: quarterly period headers, shift one quarter left, append newest quarter ---
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates, shift one left, append newest publish date ---
$ws.Range("D9").Value = "1400-10-30 (2)"
$ws.Range("E9").Value = "1401-04-08 (9)"
$ws.Range("F9").Value = "1401-04-30 (2)"
$ws.Range("G9").Value = "1401-09-15 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-01-29 (8)"
$ws.Range("J9").Value = "1401-04-30"
$ws.Range("K9").Value = "1401-09-15 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-01-29"
